# Change the table style applied to the table on slide 5 from the
# default "Medium Style 2" ({AF2116EA-DECE-4B36-BE53-6AEB78ED3CFE}) to
# "Medium Style 2 - Accent 2" ({29B43CF3-5EB2-4DAB-94D6-7D64BF7A6AA5}).
#
# Table styles can't be assigned by setting the .Style property directly
# (that throws in this host) - PowerPoint exposes table-style changes
# through Table.ApplyStyle(styleId).

$p = $ppt.ActivePresentation

$targetStyleId = "{29B43CF3-5EB2-4DAB-94D6-7D64BF7A6AA5}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($targetStyleId)
        }
    }
}
